$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-job-family"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates ---
# The "Constraint(s)" value that used to live on the Extension row (row 2, col AI)
# is removed; it now only appears on the Extension.extension row (row 4), unchanged.
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
